# Rename the "Collection" worksheet tab to "CRF" (commit: "rename Collection to CRF in tabs")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_QRS_6MWT")
$ws.Name = "CRF_QRS_6MWT"
